$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.425.65"
$ws.Range("E2").Value = "  +6.12%  "
$ws.Range("D3").Value = "1.720.98"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "'333.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "'0.3736"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.79%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "'0.3352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "'0.07366"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.02%  "
$ws.Range("D12").Value = "'1.006"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'6.346"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("D14").Value = "'20.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("E15").Value = "  +6.49%  "
$ws.Range("D16").Value = "1.724.31"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").Value = "'0.00001065"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "'0.06640"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Value = "'82.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "'16.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.13%  "
$ws.Range("D22").Value = "'6.103"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("D23").Value = "'12.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "26.400.49"
$ws.Range("E24").Value = "  +5.92%  "
$ws.Range("D25").Value = "'2.452"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'153.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("D27").Value = "'2.367"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").Value = "'1.378"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.18%  "
$ws.Range("D29").Value = "'19.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "1.922.30"
$ws.Range("E30").Value = "  +4.23%  "
$ws.Range("D31").Value = "'130.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("D32").Value = "'4.136"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "'5.911"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").Value = "'0.08598"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").Value = "'1.691"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "'12.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").Value = "'5.366"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("D39").Value = "'0.2154"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.84%  "
$ws.Range("D40").Value = "'0.06176"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").Value = "'8.478"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "'1.220"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.78%  "
$ws.Range("D43").Value = "'0.6137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").Value = "'1.005"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "'13.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "'3.901"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("E47").Value = "  +5.10%  "
$ws.Range("D48").Value = "'127.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").Value = "'2.030"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.84%  "
$ws.Range("D50").Value = "'0.07180"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").Value = "'76.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.57%  "
